$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$ws.Range("G1").Value = "Tags"
$ws.Range("G1").Select()
